# Rename the inline logo pictures in the document's headers and footers.
#
#   * BTec_Logo-Orange picture (in both headers)   : image1.jpg -> image2.jpg
#   * PearsonLogo.png picture  (in both footers)    : image2.png -> image1.png
#
# (The embedded media files themselves are untouched - this only changes
#  the display "name" recorded on each inline picture's non-visual
#  drawing properties, i.e. <wp:docPr name="...">.)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($range, [string]$newName) {
    $shape = $range.InlineShapes.Item(1)
    # Selecting first and renaming via the resulting Selection's
    # InlineShapes collection sidesteps an addressing quirk that some
    # header/footer story ranges hit when the InlineShape object is
    # renamed directly off the HeaderFooter.Range collection.
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Primary header + first-page header both hold the BTEC logo.
Rename-InlineLogo $sec.Headers.Item(1).Range "image2.jpg"
Rename-InlineLogo $sec.Headers.Item(2).Range "image2.jpg"

# Primary footer + first-page footer both hold the Pearson Edexcel logo.
Rename-InlineLogo $sec.Footers.Item(1).Range "image1.png"
Rename-InlineLogo $sec.Footers.Item(2).Range "image1.png"
